$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1977
$ws.Range("F6").Value = 4077
$ws.Range("F7").Value = 529
$ws.Range("F8").Value = 1044
$ws.Range("F12").Value = 2178
$ws.Range("F13").Value = 393
$ws.Range("F14").Value = 652561
$ws.Range("F15").Value = 1605
$ws.Range("F16").Value = 487
$ws.Range("F17").Value = 1434
$ws.Range("F18").Value = 666
$ws.Range("F20").Value = 1254
$ws.Range("F21").Value = 2185
$ws.Range("F22").Value = 1120
$ws.Range("F23").Value = 2680
$ws.Range("F24").Value = 1540
$ws.Range("F25").Value = 783
$ws.Range("F26").Value = 1521
$ws.Range("F29").Value = 1079
$ws.Range("F30").Value = 263
$ws.Range("F31").Value = 1075
$ws.Range("F34").Value = 2004
$ws.Range("F35").Value = 1347
$ws.Range("F36").Value = 565
$ws.Range("F37").Value = 1212
$ws.Range("F38").Value = 2235
$ws.Range("F39").Value = 1133
$ws.Range("F40").Value = 12
$ws.Range("F42").Value = 2559
$ws.Range("F45").Value = 3093
$ws.Range("F48").Value = 875
$ws.Range("F49").Value = 138

$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 473
$ws.Range("F11").Value = 144453
$ws.Range("F12").Value = 144453
$ws.Range("F17").Value = 94
$ws.Range("F18").Value = 224
$ws.Range("F21").Value = 407
$ws.Range("F22").Value = 407
$ws.Range("F23").Value = 117
$ws.Range("F24").Value = 79
$ws.Range("F27").Value = 534
$ws.Range("F32").Value = 328
$ws.Range("G35").Value = 71.09999999999999
$ws.Range("G36").Value = 115.2
$ws.Range("F41").Value = 86

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 3118
$ws.Range("F5").Value = 238
$ws.Range("F7").Value = 819
$ws.Range("F8").Value = 1159
$ws.Range("F9").Value = 632
$ws.Range("F10").Value = 1583
$ws.Range("F11").Value = 475
$ws.Range("F12").Value = 68
$ws.Range("F13").Value = 1834

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 819
$ws.Range("F3").Value = 1159
$ws.Range("F4").Value = 632
$ws.Range("F6").Value = 1583
$ws.Range("F7").Value = 475
$ws.Range("F9").Value = 1977
$ws.Range("F10").Value = 68
$ws.Range("F11").Value = 1834
$ws.Range("F12").Value = 4078
$ws.Range("F13").Value = 529
$ws.Range("F16").Value = 2178
$ws.Range("F17").Value = 393
$ws.Range("F18").Value = 652575
$ws.Range("F20").Value = 473
$ws.Range("F21").Value = 1605
$ws.Range("F22").Value = 144453
$ws.Range("F23").Value = 1434
$ws.Range("F24").Value = 666
$ws.Range("F26").Value = 1254
$ws.Range("F27").Value = 2185
$ws.Range("F28").Value = 1120
$ws.Range("F29").Value = 2680
$ws.Range("F30").Value = 1540
$ws.Range("F31").Value = 783
$ws.Range("F33").Value = 1521
$ws.Range("F34").Value = 407
$ws.Range("F36").Value = 117
$ws.Range("F37").Value = 1079
$ws.Range("F38").Value = 1075
$ws.Range("F40").Value = 2004
$ws.Range("F41").Value = 1347
$ws.Range("F42").Value = 1212
$ws.Range("F43").Value = 2235
$ws.Range("F44").Value = 1133
$ws.Range("F45").Value = 328
$ws.Range("F46").Value = 328
$ws.Range("F48").Value = 2559
$ws.Range("F51").Value = 3093
$ws.Range("F52").Value = 138
